$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (rich-text shared strings; all runs share identical formatting,
#     so replacing the whole cell text is visually/structurally equivalent) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Cells whose style/number-format is unchanged: just update the value ---
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 0
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 91
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = -21.551724137931
$ws.Range("L16").Value = -44.171779141104
$ws.Range("M16").Value = -24.793388429752
$ws.Range("N16").Value = -87.110481586402
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 18
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 203
$ws.Range("J17").Value = 190
$ws.Range("K17").Value = 6.842105263157
$ws.Range("L17").Value = 6.842105263157
$ws.Range("M17").Value = 73.504273504273
$ws.Range("N17").Value = 1.5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -46.666666666666
$ws.Range("I18").Value = 99
$ws.Range("J18").Value = 154
$ws.Range("K18").Value = -35.714285714285
$ws.Range("L18").Value = -31.25
$ws.Range("M18").Value = 15.116279069767
$ws.Range("N18").Value = -70.796460176991
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -30.357142857142
$ws.Range("I19").Value = 440
$ws.Range("J19").Value = 519
$ws.Range("K19").Value = -15.221579961464
$ws.Range("L19").Value = -34.131736526946
$ws.Range("M19").Value = 83.333333333333
$ws.Range("N19").Value = 10.275689223057
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = -9.259259259259
$ws.Range("L20").Value = -12.5
$ws.Range("M20").Value = 6.521739130434
$ws.Range("N20").Value = -85.714285714285
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -42.424242424242
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -29.906542056074
$ws.Range("I21").Value = 891
$ws.Range("J21").Value = 1044
$ws.Range("K21").Value = -14.655172413793
$ws.Range("L21").Value = -28.145161290322
$ws.Range("M21").Value = 43.941841680129
$ws.Range("N21").Value = -55.693684733963
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 12.5
$ws.Range("M22").Value = 12.5
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -77.777777777777
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = -5
$ws.Range("I23").Value = 185
$ws.Range("J23").Value = 149
$ws.Range("K23").Value = 24.161073825503
$ws.Range("L23").Value = 20.12987012987
$ws.Range("M23").Value = 56.779661016949
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 128
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = 7.563025210084
$ws.Range("I24").Value = 1236
$ws.Range("J24").Value = 1087
$ws.Range("K24").Value = 13.707451701931
$ws.Range("L24").Value = -37.951807228915
$ws.Range("M24").Value = 87.556904400607
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 74
$ws.Range("G25").Value = 83
$ws.Range("H25").Value = -10.843373493975
$ws.Range("I25").Value = 814
$ws.Range("J25").Value = 618
$ws.Range("K25").Value = 31.715210355987
$ws.Range("L25").Value = -49.156777014366
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 42
$ws.Range("H26").Value = 13.513513513513
$ws.Range("I26").Value = 386
$ws.Range("J26").Value = 391
$ws.Range("K26").Value = -1.278772378516
$ws.Range("L26").Value = -2.030456852791
$ws.Range("M26").Value = 38.351254480286
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -11.111111111111
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 36
$ws.Range("J28").Value = 44
$ws.Range("K28").Value = -18.181818181818
$ws.Range("L28").Value = -7.692307692307

# --- Cells converting from the text placeholder ("***.*"/blank marker) to a real number:
#     set the numeric value, then copy number-format+style from a donor cell that already
#     carries the correct integer/percent style so the cellXfs index matches exactly. ---
$ws.Range("D15").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("J14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("J14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("J14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("J14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# --- Cells converting from a number back to the text placeholder: first force a text
#     number-format so the string literal is not auto-coerced back to a number, set the
#     value, then copy format from a donor cell already on the text style. ---
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
